$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly to fit the new part numbers (stored width ends up
# rounded to the nearest pixel column width by Excel -> 14.12 "characters"
# rounds to a stored width of 15)
$ws.Columns.Item(1).ColumnWidth = 14.12

# Update the capacitor quantity (row 3, column D) from 2 to 11
$ws.Range("D3").Value = 11

# Row 11: LPC1768 microcontroller
$ws.Cells.Item(11, 1).Value = "LPC1768FBD100,551"
$ws.Cells.Item(11, 2).Formula = '=HYPERLINK("https://octopart.com/lpc1768fbd100%2C551-nxp+semiconductors-11854624","Octopart")'
$ws.Cells.Item(11, 3).Value = "uC"
$ws.Cells.Item(11, 4).Value = 1

# Row 12: 12MHz crystal
$ws.Cells.Item(12, 1).Value = "445C35A12M00000"
$ws.Cells.Item(12, 2).Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/cts-frequency-controls/445C35A12M00000/CTX1435CT-ND/5875920","Digikey - CTX1435CT-ND")'
$ws.Cells.Item(12, 3).Value = "12MHz crystal"
$ws.Cells.Item(12, 4).Value = 1

# Copy the style used for existing part-number / hyperlink cells onto the new rows
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11:A12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11:B12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Move the active selection as in the edited file
$ws.Range("C13").Select()
